# This script swaps the B:G column values between pairs of adjacent rows
# that were accidentally stored in reversed order (item code / price / qty
# / value got swapped between two rows sharing the same item description).
#
# Row pairs to swap (1-based worksheet rows, matching the XML row numbers):
#   149/150, 264/265, 279/280, 313/314, 346/347, 350/351, 355/356, 372/373,
#   379/380, 382/383, 431/432, 536/537, 581/582, 593/594, 720/721, 872/873

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(149, 150),
    @(264, 265),
    @(279, 280),
    @(313, 314),
    @(346, 347),
    @(350, 351),
    @(355, 356),
    @(372, 373),
    @(379, 380),
    @(382, 383),
    @(431, 432),
    @(536, 537),
    @(581, 582),
    @(593, 594),
    @(720, 721),
    @(872, 873)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B through G hold the item code, name, cost, price, qty and value.
    foreach ($col in @("B", "C", "D", "E", "F", "G")) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value()
        $v2 = $cell2.Value()

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

"Swapped $($rowPairs.Count) row pairs"
